$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 606.1591
$ws.Range("J17").Value = 604.8372000000001
$ws.Range("L17").Value = 1814.5116
$ws.Range("N17").Value = -2150.5116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 476.23077
$ws.Range("I55").Value = 757.2857
$ws.Range("J55").Value = 148.33333
$ws.Range("K55").Value = 757.2857
$ws.Range("L55").Value = 148.33333
$ws.Range("M55").Value = -543.2857
$ws.Range("N55").Value = -576.3333299999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1548.1111
$ws.Range("J69").Value = 1519.1923
$ws.Range("L69").Value = 4557.5769
$ws.Range("N69").Value = -6305.5769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 1548.1111
$ws.Range("J72").Value = 1519.1923
$ws.Range("L72").Value = 13672.7307
$ws.Range("N72").Value = -22408.7307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1166.6666
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4267.263
$ws.Range("I132").Value = 5251.923
$ws.Range("K132").Value = 15755.769
$ws.Range("M132").Value = -13225.769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 102148.3
$ws.Range("I137").Value = 2136.4
$ws.Range("K137").Value = 6409.200000000001
$ws.Range("M137").Value = -3859.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1300.1915
$ws.Range("I138").Value = 529.70966
$ws.Range("J138").Value = 2793
$ws.Range("K138").Value = 1589.12898
$ws.Range("L138").Value = 8379
$ws.Range("M138").Value = 3550.87102
$ws.Range("N138").Value = -18659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2493.7144
$ws.Range("I61").Value = 1668.7059
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1668.7059
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -1456.7059
$ws.Range("N61").Value = -6424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 90914910
$ws.Range("I74").Value = 142863570
$ws.Range("K74").Value = 142863570
$ws.Range("M74").Value = -142862696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 90914910
$ws.Range("I77").Value = 142863570
$ws.Range("K77").Value = 714317850
$ws.Range("M77").Value = -714313482

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15157.027
$ws.Range("I132").Value = 1518.6875
$ws.Range("K132").Value = 4556.0625
$ws.Range("M132").Value = -2026.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2493.7144
$ws.Range("I136").Value = 1668.7059
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 5006.1177
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -2456.1177
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33100.234
$ws.Range("I134").Value = 44547.48
$ws.Range("K134").Value = 133642.44
$ws.Range("M134").Value = -131107.44

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 43.076923
$ws.Range("I7").Value = 45.333332
$ws.Range("J7").Value = 38
$ws.Range("K7").Value = 45.333332
$ws.Range("L7").Value = 38
$ws.Range("M7").Value = 67.666668
$ws.Range("N7").Value = -264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19014.611
$ws.Range("I31").Value = 24404.846
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 24404.846
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -24109.846
$ws.Range("N31").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 19014.611
$ws.Range("I34").Value = 24404.846
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 24404.846
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -24202.846
$ws.Range("N34").Value = -5404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19744.408
$ws.Range("I58").Value = 1399.3529
$ws.Range("J58").Value = 50931
$ws.Range("K58").Value = 1399.3529
$ws.Range("L58").Value = 50931
$ws.Range("M58").Value = -1196.3529
$ws.Range("N58").Value = -51337

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 50767.09
$ws.Range("I132").Value = 59047.223
$ws.Range("K132").Value = 177141.669
$ws.Range("M132").Value = -174611.669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1112.091
$ws.Range("I134").Value = 927
$ws.Range("J134").Value = 1266.3334
$ws.Range("K134").Value = 2781
$ws.Range("L134").Value = 3799.0002
$ws.Range("M134").Value = -246
$ws.Range("N134").Value = -8869.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 19744.408
$ws.Range("I136").Value = 1399.3529
$ws.Range("J136").Value = 50931
$ws.Range("K136").Value = 4198.0587
$ws.Range("L136").Value = 152793
$ws.Range("M136").Value = -1648.0587
$ws.Range("N136").Value = -157893

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 794.29
$ws.Range("J131").Value = 795.24243
$ws.Range("L131").Value = 2385.72729
$ws.Range("N131").Value = -12465.72729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 108622
$ws.Range("I138").Value = 1463.8125
$ws.Range("J138").Value = 251499.58
$ws.Range("K138").Value = 4391.4375
$ws.Range("L138").Value = 754498.74
$ws.Range("M138").Value = 748.5625
$ws.Range("N138").Value = -764778.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3139.7778
$ws.Range("I97").Value = 1209.6666
$ws.Range("J97").Value = 7000
$ws.Range("K97").Value = 1209.6666
$ws.Range("L97").Value = 7000
$ws.Range("M97").Value = -713.6666
$ws.Range("N97").Value = -7992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3659.7837
$ws.Range("J126").Value = 5566.6665
$ws.Range("L126").Value = 16699.9995
$ws.Range("N126").Value = -21639.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 88191.72
$ws.Range("I132").Value = 76675.21000000001
$ws.Range("J132").Value = 128499.5
$ws.Range("K132").Value = 230025.63
$ws.Range("L132").Value = 385498.5
$ws.Range("M132").Value = -227495.63
$ws.Range("N132").Value = -390558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7033.1665
$ws.Range("I22").Value = 7000
$ws.Range("J22").Value = 7066.3335
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 7066.3335
$ws.Range("M22").Value = -6705
$ws.Range("N22").Value = -7656.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7033.1665
$ws.Range("I27").Value = 7000
$ws.Range("J27").Value = 7066.3335
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 7066.3335
$ws.Range("M27").Value = -6893
$ws.Range("N27").Value = -7280.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -67488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2329.611
$ws.Range("I132").Value = 1548.6364
$ws.Range("K132").Value = 4645.9092
$ws.Range("M132").Value = -2115.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 36886.645
$ws.Range("I136").Value = 42793.25
$ws.Range("J136").Value = 1447
$ws.Range("K136").Value = 128379.75
$ws.Range("L136").Value = 4341
$ws.Range("M136").Value = -125829.75
$ws.Range("N136").Value = -9441

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2562.125
$ws.Range("I132").Value = 1100
$ws.Range("K132").Value = 3300
$ws.Range("M132").Value = -770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 32259898
$ws.Range("J136").Value = 2110.1
$ws.Range("L136").Value = 6330.299999999999
$ws.Range("N136").Value = -11430.3
